$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2045454545454546
$ws.Range("C2").Value = 0.5426136363636364
$ws.Range("J2").Value = 0.01136363636363636
$ws.Range("P2").Value = 0.1420454545454546
$ws.Range("S2").Value = 0.09943181818181818
$ws.Range("B3").Value = 0.0202020202020202
$ws.Range("C3").Value = 0.0303030303030303
$ws.Range("J3").Value = 0.02525252525252525
$ws.Range("P3").Value = 0.7626262626262627
$ws.Range("S3").Value = 0.1616161616161616
$ws.Range("J4").Value = 0.02040816326530612
$ws.Range("P4").Value = 0.6326530612244898
$ws.Range("S4").Value = 0.3469387755102041
$ws.Range("B6").Value = 0.06374501992031872
$ws.Range("D6").Value = 0.00796812749003984
$ws.Range("F6").Value = 0.04780876494023904
$ws.Range("J6").Value = 0.3067729083665339
$ws.Range("O6").Value = 0.02788844621513944
$ws.Range("Q6").Value = 0.1593625498007968
$ws.Range("R6").Value = 0.05577689243027888
$ws.Range("S6").Value = 0.3306772908366534
$ws.Range("B7").Value = 0.09090909090909091
$ws.Range("D7").Value = 0.003636363636363636
$ws.Range("E7").Value = 0.003636363636363636
$ws.Range("F7").Value = 0.04727272727272727
$ws.Range("J7").Value = 0.1163636363636364
$ws.Range("O7").Value = 0.02545454545454546
$ws.Range("Q7").Value = 0.1272727272727273
$ws.Range("R7").Value = 0.1090909090909091
$ws.Range("S7").Value = 0.4763636363636364
$ws.Range("B8").Value = 0.1269230769230769
$ws.Range("D8").Value = 0.03076923076923077
$ws.Range("F8").Value = 0.04807692307692308
$ws.Range("J8").Value = 0.1269230769230769
$ws.Range("O8").Value = 0.03076923076923077
$ws.Range("Q8").Value = 0.1807692307692308
$ws.Range("R8").Value = 0.07884615384615384
$ws.Range("S8").Value = 0.3769230769230769
$ws.Range("B9").Value = 0.08906882591093117
$ws.Range("D9").Value = 0.004048582995951417
$ws.Range("F9").Value = 0.06072874493927125
$ws.Range("J9").Value = 0.1497975708502024
$ws.Range("O9").Value = 0.0242914979757085
$ws.Range("Q9").Value = 0.1538461538461539
$ws.Range("R9").Value = 0.06477732793522267
$ws.Range("S9").Value = 0.4534412955465587
$ws.Range("B10").Value = 0.09727081875437368
$ws.Range("D10").Value = 0.02169349195241427
$ws.Range("E10").Value = 0.0006997900629811056
$ws.Range("F10").Value = 0.07417774667599721
$ws.Range("J10").Value = 0.125962211336599
$ws.Range("O10").Value = 0.02309307207837649
$ws.Range("Q10").Value = 0.2239328201539538
$ws.Range("R10").Value = 0.08047585724282715
$ws.Range("S10").Value = 0.3526941917424772
$ws.Range("G11").Value = 0.128735632183908
$ws.Range("J11").Value = 0.0896551724137931
$ws.Range("K11").Value = 0.1954022988505747
$ws.Range("L11").Value = 0.5655172413793104
$ws.Range("S11").Value = 0.02068965517241379
$ws.Range("G12").Value = 0.80078125
$ws.Range("J12").Value = 0.14453125
$ws.Range("L12").Value = 0.0390625
$ws.Range("S12").Value = 0.015625
$ws.Range("G13").Value = 0.6666666666666666
$ws.Range("J13").Value = 0.3333333333333333
$ws.Range("F15").Value = 0.01083032490974729
$ws.Range("H15").Value = 0.1913357400722022
$ws.Range("I15").Value = 0.06498194945848375
$ws.Range("J15").Value = 0.2779783393501805
$ws.Range("K15").Value = 0.07942238267148015
$ws.Range("M15").Value = 0.003610108303249098
$ws.Range("O15").Value = 0.07220216606498195
$ws.Range("S15").Value = 0.2996389891696751
$ws.Range("F16").Value = 0.03139013452914798
$ws.Range("H16").Value = 0.1838565022421525
$ws.Range("I16").Value = 0.1210762331838565
$ws.Range("J16").Value = 0.3497757847533632
$ws.Range("K16").Value = 0.1569506726457399
$ws.Range("M16").Value = 0.01345291479820628
$ws.Range("O16").Value = 0.05381165919282511
$ws.Range("S16").Value = 0.08968609865470852
$ws.Range("F17").Value = 0.0210727969348659
$ws.Range("H17").Value = 0.1781609195402299
$ws.Range("I17").Value = 0.08237547892720307
$ws.Range("J17").Value = 0.3927203065134099
$ws.Range("K17").Value = 0.1149425287356322
$ws.Range("M17").Value = 0.01915708812260536
$ws.Range("N17").Value = 0.003831417624521073
$ws.Range("O17").Value = 0.07279693486590039
$ws.Range("S17").Value = 0.1149425287356322
$ws.Range("F18").Value = 0.01401869158878505
$ws.Range("H18").Value = 0.1869158878504673
$ws.Range("I18").Value = 0.1121495327102804
$ws.Range("J18").Value = 0.3691588785046729
$ws.Range("K18").Value = 0.1355140186915888
$ws.Range("M18").Value = 0.004672897196261682
$ws.Range("O18").Value = 0.07943925233644859
$ws.Range("S18").Value = 0.09813084112149532
$ws.Range("F19").Value = 0.02290076335877863
$ws.Range("H19").Value = 0.2061068702290076
$ws.Range("I19").Value = 0.09368494101318529
$ws.Range("J19").Value = 0.3532269257460097
$ws.Range("K19").Value = 0.1367106176266482
$ws.Range("M19").Value = 0.02081887578070784
$ws.Range("N19").Value = 0.001387925052047189
$ws.Range("O19").Value = 0.05829285218598196
$ws.Range("S19").Value = 0.1068702290076336
